# Ticket_Tally_Sheet.xlsx -- "Added values for 11/9 week."
#
# The duplicated sheet "Tally_20171019 (2)" is turned into the new
# 11/9/2017 weekly tally: it is renamed, its start date is bumped to the
# new Thursday, and the tallies for the week are filled in.

$wb = $excel.ActiveWorkbook

# --- Turn the duplicated sheet into the new 11/9 week sheet -----------
$ws = $wb.Worksheets.Item("Tally_20171019 (2)")
$ws.Name = "Tally_20171109"

# New week start date (Thursday 2017-11-09)
$ws.Range("B1").Value = 43048

# Benson Fabonan row (Assigned / Fixed)
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0

$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0

# Manuel Alberto Lomotan row (Assigned / Fixed)
$ws.Range("C9").Value = 4
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 3
$ws.Range("G9").Value = 4

$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3

# Jeric Ryan De Josef row (Assigned / Fixed)
$ws.Range("C13").Value = 4
$ws.Range("D13").Value = 3
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0

$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0

# --- Selection left on the previous (10/19) week's sheet ---------------
$wsPrev = $wb.Worksheets.Item("Tally_20171019")
$wsPrev.Activate() | Out-Null
$wsPrev.Range("C5:F6").Select() | Out-Null

# --- New sheet becomes the active / selected tab ------------------------
$ws.Activate() | Out-Null
$ws.Range("L18").Select() | Out-Null
